$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster becomes "ECs" and numeric values updated
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.619088000000001
$ws.Range("H2").Value = 13.857264
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.0552105
$ws.Range("N2").Value = 0.110421
$ws.Range("O2").Value = 0.01579120108286444
$ws.Range("P2").Value = 0.01113815545262342
$ws.Range("Q2").Value = 0.255022158024
$ws.Range("R2").Value = 1.530132948144
$ws.Range("S2").Value = 0.01579120108286444
$ws.Range("T2").Value = 0.01113815545262342

# Row 3: Target cluster becomes "FAPs" and numeric values updated
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.619088000000001
$ws.Range("H3").Value = 13.857264
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.655685
$ws.Range("N3").Value = 7.967055
$ws.Range("O3").Value = 0.7595739188695422
$ws.Range("P3").Value = 0.8036360573586609
$ws.Range("Q3").Value = 12.26684271528
$ws.Range("R3").Value = 110.40158443752
$ws.Range("S3").Value = 0.7595739188695422
$ws.Range("T3").Value = 0.8036360573586609

# Row 4 (new): Target cluster "Neutro"
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.619088000000001
$ws.Range("H4").Value = 13.857264
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.26551
$ws.Range("N4").Value = 0.79653
$ws.Range("O4").Value = 0.07594065982940451
$ws.Range("P4").Value = 0.08034590306805892
$ws.Range("Q4").Value = 1.22641405488
$ws.Range("R4").Value = 11.03772649392
$ws.Range("S4").Value = 0.07594065982940451
$ws.Range("T4").Value = 0.08034590306805892

# Row 5 (new): Target cluster "sCs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.619088000000001
$ws.Range("H5").Value = 13.857264
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.519877
$ws.Range("N5").Value = 1.039754
$ws.Range("O5").Value = 0.1486942202181889
$ws.Range("P5").Value = 0.1048798841206565
$ws.Range("Q5").Value = 2.401357612176001
$ws.Range("R5").Value = 14.408145673056
$ws.Range("S5").Value = 0.1486942202181889
$ws.Range("T5").Value = 0.1048798841206565
